# //Changed to pickup store
# Add new automation test case rows to the "Test Cases" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

$ws.Cells.Item(26, 1).Value = "Mandatory fields validation on review page"
$ws.Cells.Item(26, 2).Value = "Done"

$ws.Cells.Item(27, 1).Value = "Verify Everyday Value Offer text"
$ws.Cells.Item(27, 2).Value = "Done"

$ws.Cells.Item(28, 1).Value = "Verify alert in deal of 2nd step "
$ws.Cells.Item(28, 2).Value = "Done"

$ws.Cells.Item(29, 1).Value = "Verify alert if select no pizza"
$ws.Cells.Item(29, 2).Value = "Done"

# Make the "nd" in "2nd" a superscript to match the ordinal formatting
$run = $ws.Cells.Item(28, 1).Characters(26, 2)
$run.Font.Superscript = $true

# The superscript run makes row 28 slightly taller than the default row height
$ws.Rows.Item(28).RowHeight = 13.4

$ws.Range("A30").Select()
